# Remove the "pandasgui" instructions block:
#   - "To show output on gui we are using pandasgui so install it using "
#   - ">pip install pandasgui"
#   - (empty italic ListParagraph)
#   - (empty italic paragraph, ind left=720)
# while leaving the following empty paragraph (plain <w:rPr><w:i/></w:rPr>)
# and everything else untouched.

$d = $word.ActiveDocument

# Locate the first paragraph of the block to remove via text search so the
# script is resilient to any paragraph-count drift elsewhere in the doc.
# NB: capture the Range once and re-use it -- Find.Execute mutates the
# Range object it is called on in place, so a fresh $d.Content afterwards
# would still point at the whole document.
$rng = $d.Content
$found = $rng.Find.Execute("To show output on", $true, $false, $false, `
    $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not locate the 'To show output on' paragraph"
}

$startPara = $rng.Paragraphs(1)

# Walk forward three more paragraphs (the ">pip install pandasgui" line and
# the two blank spacer paragraphs) so $endPara lands on the paragraph that
# must be preserved.
$endPara = $startPara
for ($i = 0; $i -lt 3; $i++) {
    $endPara = $endPara.Next()
}
$endPara = $endPara.Next()

$killRange = $d.Range($startPara.Range.Start, $endPara.Range.Start)
$killRange.Delete()

Write-Host "Removed pandasgui block; paragraphs now: $($d.Paragraphs.Count)"
